$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1847.5
$ws.Cells.Item(40, 9).Value = 1935.1852
$ws.Cells.Item(40, 10).Value = 1708.2354
$ws.Cells.Item(40, 11).Value = 1935.1852
$ws.Cells.Item(40, 12).Value = 1708.2354
$ws.Cells.Item(40, 13).Value = -1760.1852
$ws.Cells.Item(40, 14).Value = -2058.2354

$ws.Cells.Item(53, 8).Value = 277.5
$ws.Cells.Item(53, 9).Value = 161
$ws.Cells.Item(53, 11).Value = 161
$ws.Cells.Item(53, 13).Value = 476

$ws.Cells.Item(58, 8).Value = 1780.6364
$ws.Cells.Item(58, 9).Value = 176.33333
$ws.Cells.Item(58, 10).Value = 9000
$ws.Cells.Item(58, 11).Value = 528.99999
$ws.Cells.Item(58, 12).Value = 27000
$ws.Cells.Item(58, 13).Value = -378.99999
$ws.Cells.Item(58, 14).Value = -27300

$ws.Cells.Item(62, 8).Value = 26317072
$ws.Cells.Item(62, 9).Value = 963.4167
$ws.Cells.Item(62, 11).Value = 963.4167
$ws.Cells.Item(62, 13).Value = -339.4167

$ws.Cells.Item(64, 8).Value = 2872.289
$ws.Cells.Item(64, 9).Value = 2797.8723
$ws.Cells.Item(64, 10).Value = 2969.4443
$ws.Cells.Item(64, 11).Value = 2797.8723
$ws.Cells.Item(64, 12).Value = 2969.4443
$ws.Cells.Item(64, 13).Value = -2549.8723
$ws.Cells.Item(64, 14).Value = -3465.4443

$ws.Cells.Item(65, 8).Value = 26317072
$ws.Cells.Item(65, 9).Value = 963.4167
$ws.Cells.Item(65, 11).Value = 4817.0835
$ws.Cells.Item(65, 13).Value = -1697.0835

$ws.Cells.Item(67, 8).Value = 2872.289
$ws.Cells.Item(67, 9).Value = 2797.8723
$ws.Cells.Item(67, 10).Value = 2969.4443
$ws.Cells.Item(67, 11).Value = 2797.8723
$ws.Cells.Item(67, 12).Value = 2969.4443
$ws.Cells.Item(67, 13).Value = -1939.8723
$ws.Cells.Item(67, 14).Value = -4685.4443

$ws.Cells.Item(107, 8).Value = 2200.25
$ws.Cells.Item(107, 9).Value = 2236.6365
$ws.Cells.Item(107, 10).Value = 1800
$ws.Cells.Item(107, 11).Value = 2236.6365
$ws.Cells.Item(107, 12).Value = 1800
$ws.Cells.Item(107, 13).Value = -316.6365000000001
$ws.Cells.Item(107, 14).Value = -5640

$ws.Cells.Item(132, 8).Value = 2523.4443
$ws.Cells.Item(132, 9).Value = 1715.4524
$ws.Cells.Item(132, 11).Value = 5146.357199999999
$ws.Cells.Item(132, 13).Value = -2616.357199999999

$ws.Cells.Item(138, 8).Value = 2392.2258
$ws.Cells.Item(138, 9).Value = 2841.3333
$ws.Cells.Item(138, 10).Value = 2284.44
$ws.Cells.Item(138, 11).Value = 8523.999899999999
$ws.Cells.Item(138, 12).Value = 6853.32
$ws.Cells.Item(138, 13).Value = -3383.999899999999
$ws.Cells.Item(138, 14).Value = -17133.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1474.4783
$ws.Cells.Item(2, 9).Value = 952.3611
$ws.Cells.Item(2, 10).Value = 3354.1
$ws.Cells.Item(2, 11).Value = 952.3611
$ws.Cells.Item(2, 12).Value = 3354.1
$ws.Cells.Item(2, 13).Value = -839.3611
$ws.Cells.Item(2, 14).Value = -3580.1

$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).ClearContents()
$ws.Cells.Item(13, 13).ClearContents()
$ws.Cells.Item(13, 14).Value = 0

$ws.Cells.Item(74, 8).Value = 1434.7646
$ws.Cells.Item(74, 9).Value = 767.4838999999999
$ws.Cells.Item(74, 10).Value = 8330
$ws.Cells.Item(74, 11).Value = 767.4838999999999
$ws.Cells.Item(74, 12).Value = 8330
$ws.Cells.Item(74, 13).Value = 106.5161000000001
$ws.Cells.Item(74, 14).Value = -10078

$ws.Cells.Item(77, 8).Value = 1434.7646
$ws.Cells.Item(77, 9).Value = 767.4838999999999
$ws.Cells.Item(77, 10).Value = 8330
$ws.Cells.Item(77, 11).Value = 3837.4195
$ws.Cells.Item(77, 12).Value = 41650
$ws.Cells.Item(77, 13).Value = 530.5805
$ws.Cells.Item(77, 14).Value = -50386

$ws.Cells.Item(97, 8).Value = 743.8461
$ws.Cells.Item(97, 9).Value = 145.42857
$ws.Cells.Item(97, 11).Value = 145.42857
$ws.Cells.Item(97, 13).Value = 350.57143

$ws.Cells.Item(116, 8).Value = 1474.4783
$ws.Cells.Item(116, 9).Value = 952.3611
$ws.Cells.Item(116, 10).Value = 3354.1
$ws.Cells.Item(116, 11).Value = 952.3611
$ws.Cells.Item(116, 12).Value = 3354.1
$ws.Cells.Item(116, 13).Value = 1341.6389
$ws.Cells.Item(116, 14).Value = -7942.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1474.4783
$ws.Cells.Item(3, 9).Value = 952.3611
$ws.Cells.Item(3, 10).Value = 3354.1
$ws.Cells.Item(3, 11).Value = 952.3611
$ws.Cells.Item(3, 12).Value = 3354.1
$ws.Cells.Item(3, 13).Value = -838.3611
$ws.Cells.Item(3, 14).Value = -3582.1

$ws.Cells.Item(10, 8).Value = 5005
$ws.Cells.Item(10, 9).Value = 5005
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 5005
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).ClearContents()
$ws.Cells.Item(10, 14).Value = -4865

$ws.Cells.Item(126, 8).Value = 23935
$ws.Cells.Item(126, 10).Value = 23935
$ws.Cells.Item(126, 12).Value = 23935
$ws.Cells.Item(126, 14).Value = -33815

$ws.Cells.Item(134, 8).Value = 1802.6923
$ws.Cells.Item(134, 9).Value = 1786.25
$ws.Cells.Item(134, 11).Value = 5358.75
$ws.Cells.Item(134, 13).Value = -2823.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 9106
$ws.Cells.Item(50, 10).Value = 9106
$ws.Cells.Item(50, 12).Value = 9106
$ws.Cells.Item(50, 14).Value = -10356

$ws.Cells.Item(99, 8).Value = 1853.6897
$ws.Cells.Item(99, 9).Value = 1675.5
$ws.Cells.Item(99, 10).Value = 2020
$ws.Cells.Item(99, 11).Value = 1675.5
$ws.Cells.Item(99, 12).Value = 2020
$ws.Cells.Item(99, 13).Value = -177.5
$ws.Cells.Item(99, 14).Value = -5016

$ws.Cells.Item(109, 8).Value = 27800
$ws.Cells.Item(109, 10).Value = 27800
$ws.Cells.Item(109, 12).Value = 27800
$ws.Cells.Item(109, 14).Value = -29880

$ws.Cells.Item(126, 8).Value = 1853.6897
$ws.Cells.Item(126, 9).Value = 1675.5
$ws.Cells.Item(126, 10).Value = 2020
$ws.Cells.Item(126, 11).Value = 5026.5
$ws.Cells.Item(126, 12).Value = 6060
$ws.Cells.Item(126, 13).Value = -2556.5
$ws.Cells.Item(126, 14).Value = -11000

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 2274768.2
$ws.Cells.Item(2, 9).Value = 100040.9
$ws.Cells.Item(2, 11).Value = 600245.3999999999
$ws.Cells.Item(2, 13).Value = -600132.3999999999

$ws.Cells.Item(112, 8).Value = 33374964
$ws.Cells.Item(112, 9).Value = 913.5
$ws.Cells.Item(112, 10).Value = 36277056
$ws.Cells.Item(112, 11).Value = 2740.5
$ws.Cells.Item(112, 12).Value = 108831168
$ws.Cells.Item(112, 13).Value = -1632.5
$ws.Cells.Item(112, 14).Value = -108833384

$ws.Cells.Item(113, 8).Value = 789.2
$ws.Cells.Item(113, 9).Value = 534.5
$ws.Cells.Item(113, 10).Value = 881.8182
$ws.Cells.Item(113, 11).Value = 1603.5
$ws.Cells.Item(113, 12).Value = 2645.4546
$ws.Cells.Item(113, 13).Value = 566.5
$ws.Cells.Item(113, 14).Value = -6985.4546

$ws.Cells.Item(115, 8).Value = 1857.1428
$ws.Cells.Item(115, 9).Value = 1280
$ws.Cells.Item(115, 11).Value = 3840
$ws.Cells.Item(115, 13).Value = -2665

$ws.Cells.Item(118, 8).Value = 1537.7778
$ws.Cells.Item(118, 9).Value = 490
$ws.Cells.Item(118, 10).Value = 3633.3333
$ws.Cells.Item(118, 11).Value = 1470
$ws.Cells.Item(118, 12).Value = 10899.9999
$ws.Cells.Item(118, 13).Value = -227
$ws.Cells.Item(118, 14).Value = -13385.9999

$ws.Cells.Item(121, 8).Value = 513825.38
$ws.Cells.Item(121, 9).Value = 333.07693
$ws.Cells.Item(121, 10).Value = 1027317.7
$ws.Cells.Item(121, 11).Value = 999.2307900000001
$ws.Cells.Item(121, 12).Value = 3081953.1
$ws.Cells.Item(121, 13).Value = 310.7692099999999
$ws.Cells.Item(121, 14).Value = -3084573.1

$ws.Cells.Item(125, 8).Value = 1946.1111

$ws.Cells.Item(126, 8).Value = 1726.5
$ws.Cells.Item(126, 10).Value = 1763.1578
$ws.Cells.Item(126, 12).Value = 5289.4734
$ws.Cells.Item(126, 14).Value = -15169.4734

$ws.Cells.Item(127, 8).Value = 631.1
$ws.Cells.Item(127, 10).Value = 631.1
$ws.Cells.Item(127, 12).Value = 1893.3
$ws.Cells.Item(127, 14).Value = -11813.3

$ws.Cells.Item(131, 8).Value = 906.3196
$ws.Cells.Item(131, 9).Value = 601.53845
$ws.Cells.Item(131, 10).Value = 953.4881
$ws.Cells.Item(131, 11).Value = 1804.61535
$ws.Cells.Item(131, 12).Value = 2860.4643
$ws.Cells.Item(131, 13).Value = 3235.38465
$ws.Cells.Item(131, 14).Value = -12940.4643

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 17537.1
$ws.Cells.Item(57, 10).Value = 18924
$ws.Cells.Item(57, 12).Value = 18924
$ws.Cells.Item(57, 14).Value = -20564

$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).ClearContents()
$ws.Cells.Item(121, 14).Value = 0

$ws.Cells.Item(135, 8).Value = 49729.566
$ws.Cells.Item(135, 10).Value = 49729.566
$ws.Cells.Item(135, 12).Value = 49729.566
$ws.Cells.Item(135, 14).Value = -59869.566

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).ClearContents()
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = 0

$ws.Cells.Item(22, 8).Value = 1222.2222
$ws.Cells.Item(22, 9).Value = 914.2857
$ws.Cells.Item(22, 10).Value = 1418.1818
$ws.Cells.Item(22, 11).Value = 914.2857
$ws.Cells.Item(22, 12).Value = 1418.1818
$ws.Cells.Item(22, 13).Value = -619.2857
$ws.Cells.Item(22, 14).Value = -2008.1818

$ws.Cells.Item(27, 8).Value = 1222.2222
$ws.Cells.Item(27, 9).Value = 914.2857
$ws.Cells.Item(27, 10).Value = 1418.1818
$ws.Cells.Item(27, 11).Value = 914.2857
$ws.Cells.Item(27, 12).Value = 1418.1818
$ws.Cells.Item(27, 13).Value = -807.2857
$ws.Cells.Item(27, 14).Value = -1632.1818

$ws.Cells.Item(46, 8).Value = 1791.5769
$ws.Cells.Item(46, 9).Value = 1655.6111
$ws.Cells.Item(46, 10).Value = 2097.5
$ws.Cells.Item(46, 11).Value = 1655.6111
$ws.Cells.Item(46, 12).Value = 2097.5
$ws.Cells.Item(46, 13).Value = -1467.6111
$ws.Cells.Item(46, 14).Value = -2473.5

$ws.Cells.Item(55, 8).Value = 597
$ws.Cells.Item(55, 9).Value = 240.54546
$ws.Cells.Item(55, 10).Value = 1157.1428
$ws.Cells.Item(55, 11).Value = 240.54546
$ws.Cells.Item(55, 12).Value = 1157.1428
$ws.Cells.Item(55, 13).Value = -67.54545999999999
$ws.Cells.Item(55, 14).Value = -1503.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 18288.5
$ws.Cells.Item(109, 10).Value = 18288.5
$ws.Cells.Item(109, 12).Value = 18288.5
$ws.Cells.Item(109, 14).Value = -21062.5
